$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B1: value 0, bold font + thin box border, centered horizontal / top vertical alignment
$ws.Range("B1").Value = 0
$ws.Range("B1").Font.Bold = $true
$ws.Range("B1").HorizontalAlignment = -4108  # xlCenter
$ws.Range("B1").VerticalAlignment = -4160    # xlTop
$ws.Range("B1").Borders.LineStyle = 1        # xlContinuous
$ws.Range("B1").Borders.Weight = 2           # xlThin

# A2: value 0, same formatting as B1 (copy the format so it reuses the same style)
$ws.Range("A2").Value = 0
$ws.Range("B1").Copy()
$ws.Range("A2").PasteSpecial(-4122)  # xlPasteFormats

# B2: plain text label (becomes a shared string), no special style
$ws.Range("B2").Value = "disconnected_elements"
